$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Duplicate the last date-pair of columns (DF:DG = "07-16_A"/"07-16_0")
#    into two brand-new columns (DH:DI), matching the next date ("07-17").
#    Copy values first, then formats, since a single combined paste loses styles.
$ws.Range("DF2:DG179").Copy() | Out-Null
$ws.Range("DH2").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("DF2:DG179").Copy() | Out-Null
$ws.Range("DH2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 2. Header row: new headers continuing the date series, same style as DG1/DF1.
$ws.Range("DF1").Copy() | Out-Null
$ws.Range("DH1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("DG1").Copy() | Out-Null
$ws.Range("DI1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("DH1").Value = "07-17_A"
$ws.Range("DI1").Value = "07-17_0"

# 3. The original last column (DG) held its running total as text (inline string
#    that merely looked numeric). Now that DI carries that text copy onward,
#    convert DG's own stored values to real numbers for every populated row.
for ($r = 2; $r -le 179; $r++) {
    $cell = $ws.Cells.Item($r, 111)  # column DG
    $t = $cell.Text
    if ($t -ne "") {
        $cell.Value = [double]$t
    }
}
